$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.818.57"
$ws.Range("E2").Value = "'  -2.29%  "
$ws.Range("D3").Value = "'3.114.36"
$ws.Range("E3").Value = "'  -0.61%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'590.44"
$ws.Range("E5").Value = "'  -2.27%  "
$ws.Range("D6").Value = "'135.76"
$ws.Range("E6").Value = "'  -4.98%  "
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'3.111.93"
$ws.Range("E8").Value = "'  -0.57%  "
$ws.Range("D9").Value = "'0.513"
$ws.Range("E9").Value = "'  -1.70%  "
$ws.Range("E10").Value = "'  -4.06%  "
$ws.Range("D11").Value = "'5.21"
$ws.Range("E11").Value = "'  -3.19%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "'  -3.12%  "
$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "'  -5.07%  "
$ws.Range("D14").Value = "'33.88"
$ws.Range("E14").Value = "'  -3.67%  "
$ws.Range("D15").Value = "'3.631.61"
$ws.Range("E15").Value = "'  -0.59%  "
$ws.Range("E16").Value = "'  +1.53%  "
$ws.Range("D17").Value = "'62.968.51"
$ws.Range("E17").Value = "'  -2.15%  "
$ws.Range("D18").Value = "'3.110.24"
$ws.Range("E18").Value = "'  -0.57%  "
$ws.Range("D19").Value = "'6.62"
$ws.Range("E19").Value = "'  -3.36%  "
$ws.Range("D20").Value = "'468.84"
$ws.Range("E20").Value = "'  -2.03%  "
$ws.Range("D21").Value = "'14.02"
$ws.Range("E21").Value = "'  -3.60%  "
$ws.Range("D22").Value = "'0.692"
$ws.Range("E22").Value = "'  -2.46%  "
$ws.Range("D23").Value = "'7.61"
$ws.Range("E23").Value = "'  -0.66%  "
$ws.Range("D24").Value = "'85.34"
$ws.Range("E24").Value = "'  +0.33%  "
$ws.Range("D25").Value = "'12.84"
$ws.Range("E25").Value = "'  -3.96%  "
$ws.Range("E26").Value = "'  +0.10%  "
$ws.Range("E27").Value = "'  -1.57%  "
$ws.Range("D28").Value = "'7.79"
$ws.Range("E28").Value = "'  -6.93%  "
$ws.Range("D29").Value = "'2.07"
$ws.Range("E29").Value = "'  +1.70%  "
$ws.Range("E30").Value = "'  -5.79%  "
$ws.Range("E31").Value = "'  -0.01%  "
$ws.Range("D32").Value = "'26.48"
$ws.Range("E32").Value = "'  -1.49%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "'  -5.64%  "
$ws.Range("D34").Value = "'2.50"
$ws.Range("E34").Value = "'  -5.20%  "
$ws.Range("E35").Value = "'  -3.11%  "
$ws.Range("D36").Value = "'5.72"
$ws.Range("E36").Value = "'  -3.87%  "
$ws.Range("D37").Value = "'51.84"
$ws.Range("E37").Value = "'  -1.04%  "
$ws.Range("D38").Value = "'0.0₃0686"
$ws.Range("E38").Value = "'  -11.06%  "
$ws.Range("E39").Value = "'  -2.09%  "
$ws.Range("D40").Value = "'414.63"
$ws.Range("E40").Value = "'  -6.68%  "
$ws.Range("D41").Value = "'8.16"
$ws.Range("E41").Value = "'  -0.27%  "
$ws.Range("D42").Value = "'2.893.78"
$ws.Range("E42").Value = "'  +1.46%  "
$ws.Range("E43").Value = "'  -12.41%  "
$ws.Range("E44").Value = "'  -5.74%  "
$ws.Range("D45").Value = "'0.259"
$ws.Range("E45").Value = "'  -0.12%  "
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "'  -6.21%  "
$ws.Range("D48").Value = "'25.24"
$ws.Range("E48").Value = "'  -2.74%  "
$ws.Range("E49").Value = "'  -0.60%  "
$ws.Range("E50").Value = "'  -7.83%  "
$ws.Range("D51").Value = "'120.06"
$ws.Range("E51").Value = "'  +0.32%  "
